$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 43
$ws.Range("H43").Value = 14878.571
$ws.Range("I43").Value = 50250.5
$ws.Range("J43").Value = 729.8
$ws.Range("K43").Value = 50250.5
$ws.Range("L43").Value = 729.8
$ws.Range("M43").Value = -50181.5
$ws.Range("N43").Value = -867.8
# Row 51
$ws.Range("H51").Value = 27187.75
$ws.Range("I51").Value = 51875.5
$ws.Range("J51").Value = 2500
$ws.Range("K51").Value = 51875.5
$ws.Range("L51").Value = 2500
$ws.Range("M51").Value = -51391.5
$ws.Range("N51").Value = -3468
# Row 62
$ws.Range("H62").Value = 1834
$ws.Range("J62").Value = 2224.75
$ws.Range("L62").Value = 2224.75
$ws.Range("N62").Value = -3472.75
# Row 65
$ws.Range("H65").Value = 1834
$ws.Range("J65").Value = 2224.75
$ws.Range("L65").Value = 11123.75
$ws.Range("N65").Value = -17363.75
# Row 80
$ws.Range("H80").Value = 4654.6
$ws.Range("I80").Value = 167.06667
$ws.Range("J80").Value = 11385.9
$ws.Range("K80").Value = 501.20001
$ws.Range("L80").Value = 34157.7
$ws.Range("M80").Value = 496.79999
$ws.Range("N80").Value = -36153.7
# Row 83
$ws.Range("H83").Value = 4654.6
$ws.Range("I83").Value = 167.06667
$ws.Range("J83").Value = 11385.9
$ws.Range("K83").Value = 1503.60003
$ws.Range("L83").Value = 102473.1
$ws.Range("M83").Value = 3488.39997
$ws.Range("N83").Value = -112457.1
# Row 127
$ws.Range("H127").Value = 1236.5333
$ws.Range("I127").Value = 331.85715
$ws.Range("J127").Value = 2028.125
$ws.Range("K127").Value = 995.5714499999999
$ws.Range("L127").Value = 6084.375
$ws.Range("M127").Value = 3964.42855
$ws.Range("N127").Value = -16004.375
# Row 129
$ws.Range("H129").Value = 1096.9166
$ws.Range("I129").Value = 452.2
$ws.Range("J129").Value = 1145.0299
$ws.Range("K129").Value = 1356.6
$ws.Range("L129").Value = 3435.0897
$ws.Range("M129").Value = 3643.4
$ws.Range("N129").Value = -13435.0897
# Row 137
$ws.Range("H137").Value = 2106.2341
$ws.Range("I137").Value = 1395.381
$ws.Range("J137").Value = 2680.3845
$ws.Range("K137").Value = 4186.143
$ws.Range("L137").Value = 8041.1535
$ws.Range("M137").Value = -1636.143
$ws.Range("N137").Value = -13141.1535
# Row 138
$ws.Range("H138").Value = 3598.3406
$ws.Range("I138").Value = 2083.0386
$ws.Range("J138").Value = 4204.4614
$ws.Range("K138").Value = 6249.1158
$ws.Range("L138").Value = 12613.3842
$ws.Range("M138").Value = -1109.1158
$ws.Range("N138").Value = -22893.3842

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 10263.962
$ws.Range("I32").Value = 9923.137000000001
$ws.Range("K32").Value = 9923.137000000001
$ws.Range("M32").Value = -9636.137000000001
# Row 122
$ws.Range("H122").Value = 7211.478
$ws.Range("I122").Value = 7448.0557
$ws.Range("K122").Value = 22344.1671
$ws.Range("M122").Value = -19894.1671
# Row 124
$ws.Range("H124").Value = 34776
$ws.Range("J124").Value = 34776
$ws.Range("L124").Value = 34776
$ws.Range("N124").Value = -44596
# Row 131
$ws.Range("H131").Value = 42586.91
$ws.Range("J131").Value = 42586.91
$ws.Range("L131").Value = 42586.91
$ws.Range("N131").Value = -52666.91
# Row 132
$ws.Range("H132").Value = 4030.25
$ws.Range("I132").Value = 4196.41
$ws.Range("J132").Value = 3531.7693
$ws.Range("K132").Value = 12589.23
$ws.Range("L132").Value = 10595.3079
$ws.Range("M132").Value = -10059.23
$ws.Range("N132").Value = -15655.3079

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 937.3125
$ws.Range("I94").Value = 909.6667
$ws.Range("J94").Value = 972.8570999999999
$ws.Range("K94").Value = 909.6667
$ws.Range("L94").Value = 972.8570999999999
$ws.Range("M94").Value = -458.6667
$ws.Range("N94").Value = -1874.8571
# Row 107
$ws.Range("H107").Value = 25574.738
$ws.Range("I107").Value = 33012.94
$ws.Range("K107").Value = 33012.94
$ws.Range("M107").Value = -31092.94

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 23
$ws.Range("H23").Value = 37752.25
$ws.Range("I23").Value = 33669.668
$ws.Range("J23").Value = 50000
$ws.Range("K23").Value = 33669.668
$ws.Range("L23").Value = 50000
$ws.Range("M23").Value = -33429.668
$ws.Range("N23").Value = -50480
# Row 27
$ws.Range("H27").Value = 37752.25
$ws.Range("I27").Value = 33669.668
$ws.Range("J27").Value = 50000
$ws.Range("K27").Value = 33669.668
$ws.Range("L27").Value = 50000
$ws.Range("M27").Value = -33477.668
$ws.Range("N27").Value = -50384
# Row 31
$ws.Range("H31").Value = 1792.8494
$ws.Range("I31").Value = 2195.25
$ws.Range("J31").Value = 1478.7805
$ws.Range("K31").Value = 2195.25
$ws.Range("L31").Value = 1478.7805
$ws.Range("M31").Value = -1900.25
$ws.Range("N31").Value = -2068.7805
# Row 34
$ws.Range("H34").Value = 1792.8494
$ws.Range("I34").Value = 2195.25
$ws.Range("J34").Value = 1478.7805
$ws.Range("K34").Value = 2195.25
$ws.Range("L34").Value = 1478.7805
$ws.Range("M34").Value = -1993.25
$ws.Range("N34").Value = -1882.7805
# Row 99
$ws.Range("H99").Value = 2042.8334
$ws.Range("I99").Value = 2000
$ws.Range("J99").Value = 2171.3333
$ws.Range("K99").Value = 2000
$ws.Range("L99").Value = 2171.3333
$ws.Range("M99").Value = -502
$ws.Range("N99").Value = -5167.3333
# Row 126
$ws.Range("H126").Value = 2042.8334
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 2171.3333
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 6513.999899999999
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -11453.9999

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 4
$ws.Range("H4").Value = 50
$ws.Range("I4").Value = 50
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 150
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -38
# Row 5
$ws.Range("H5").Value = 848.70215
$ws.Range("I5").Value = 1137.0741
$ws.Range("J5").Value = 459.4
$ws.Range("K5").Value = 3411.2223
$ws.Range("L5").Value = 1378.2
$ws.Range("M5").Value = -3299.2223
$ws.Range("N5").Value = -1602.2
# Row 68
$ws.Range("H68").Value = 1157.8695
$ws.Range("I68").Value = 834.11536
$ws.Range("J68").Value = 1578.75
$ws.Range("K68").Value = 2502.34608
$ws.Range("L68").Value = 4736.25
$ws.Range("M68").Value = -1691.34608
$ws.Range("N68").Value = -6358.25
# Row 69
$ws.Range("H69").Value = 1245
$ws.Range("I69").Value = 993.3333
$ws.Range("K69").Value = 2979.9999
$ws.Range("M69").Value = -2168.9999
# Row 71
$ws.Range("H71").Value = 1157.8695
$ws.Range("I71").Value = 834.11536
$ws.Range("J71").Value = 1578.75
$ws.Range("K71").Value = 7507.03824
$ws.Range("L71").Value = 14208.75
$ws.Range("M71").Value = -3451.03824
$ws.Range("N71").Value = -22320.75
# Row 72
$ws.Range("H72").Value = 1245
$ws.Range("I72").Value = 993.3333
$ws.Range("K72").Value = 8939.9997
$ws.Range("M72").Value = -4883.9997
# Row 104
$ws.Range("H104").Value = 5998.75
$ws.Range("J104").Value = 5998.75
$ws.Range("L104").Value = 17996.25
$ws.Range("N104").Value = -23238.25
# Row 107
$ws.Range("H107").Value = 1273.5933
$ws.Range("I107").Value = 1223.2222
$ws.Range("J107").Value = 1352.4348
$ws.Range("K107").Value = 3669.6666
$ws.Range("L107").Value = 4057.3044
$ws.Range("M107").Value = -1749.6666
$ws.Range("N107").Value = -7897.3044
# Row 119
$ws.Range("H119").Value = 7363.5454
$ws.Range("I119").Value = 3799.8
$ws.Range("K119").Value = 11399.4
$ws.Range("M119").Value = -6561.400000000001
# Row 120
$ws.Range("H120").Value = 9706
$ws.Range("I120").Value = 7132.5
$ws.Range("K120").Value = 21397.5
$ws.Range("M120").Value = -16559.5
# Row 131
$ws.Range("H131").Value = 3828.639
$ws.Range("I131").Value = 504.58334
$ws.Range("J131").Value = 5490.6665
$ws.Range("K131").Value = 1513.75002
$ws.Range("L131").Value = 16471.9995
$ws.Range("M131").Value = 3526.24998
$ws.Range("N131").Value = -26551.9995
# Row 135
$ws.Range("H135").Value = 848.70215
$ws.Range("I135").Value = 1137.0741
$ws.Range("J135").Value = 459.4
$ws.Range("K135").Value = 10233.6669
$ws.Range("L135").Value = 4134.599999999999
$ws.Range("M135").Value = -7698.6669
$ws.Range("N135").Value = -9204.599999999999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 5
$ws.Range("H5").Value = 32778.816
$ws.Range("I5").Value = 3888
$ws.Range("J5").Value = 49287.855
$ws.Range("K5").Value = 3888
$ws.Range("L5").Value = 49287.855
$ws.Range("M5").Value = -3776
$ws.Range("N5").Value = -49511.855
# Row 102
$ws.Range("H102").Value = 3359.5667
$ws.Range("I102").Value = 3362.5186
$ws.Range("J102").Value = 3333
$ws.Range("K102").Value = 3362.5186
$ws.Range("L102").Value = 3333
$ws.Range("M102").Value = -1740.5186
$ws.Range("N102").Value = -6577
# Row 103
$ws.Range("H103").Value = 22666.666
$ws.Range("J103").Value = 22666.666
$ws.Range("L103").Value = 22666.666
$ws.Range("N103").Value = -25010.666
# Row 132
$ws.Range("H132").Value = 2527.7334
$ws.Range("I132").Value = 2459.6667
$ws.Range("J132").Value = 2800
$ws.Range("K132").Value = 7379.000100000001
$ws.Range("L132").Value = 8400
$ws.Range("M132").Value = -4849.000100000001
$ws.Range("N132").Value = -13460

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 20458342
$ws.Range("I122").Value = 22729772
$ws.Range("J122").Value = 18186910
$ws.Range("K122").Value = 68189316
$ws.Range("L122").Value = 54560730
$ws.Range("M122").Value = -68186866
$ws.Range("N122").Value = -54565630
# Row 132
$ws.Range("H132").Value = 5050.263
$ws.Range("I132").Value = 5136.9443
$ws.Range("J132").Value = 3490
$ws.Range("K132").Value = 15410.8329
$ws.Range("L132").Value = 10470
$ws.Range("M132").Value = -12880.8329
$ws.Range("N132").Value = -15530

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 34
$ws.Range("H34").Value = 45000
$ws.Range("J34").Value = 40000
$ws.Range("L34").Value = 40000
$ws.Range("N34").Value = -40406
# Row 126
$ws.Range("H126").Value = 4265.7144
$ws.Range("I126").Value = 4643.6665
$ws.Range("J126").Value = 1998
$ws.Range("K126").Value = 13930.9995
$ws.Range("L126").Value = 5994
$ws.Range("M126").Value = -11460.9995
$ws.Range("N126").Value = -10934
# Row 136
$ws.Range("H136").Value = 2534.6858
$ws.Range("I136").Value = 2200.3914
$ws.Range("J136").Value = 3175.4167
$ws.Range("K136").Value = 6601.174199999999
$ws.Range("L136").Value = 9526.250100000001
$ws.Range("M136").Value = -4051.174199999999
$ws.Range("N136").Value = -14626.2501
